$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3757806455830632
$ws.Cells.Item(2, 3).Value = 0.05882550405871712
$ws.Cells.Item(2, 4).Value = 0.6536584357451716
$ws.Cells.Item(2, 5).Value = 0.2667908287872081
$ws.Cells.Item(2, 7).Value = 0.7580791313857489
$ws.Cells.Item(2, 8).Value = 0.8544282641165424
$ws.Cells.Item(2, 9).Value = 0.6717967527746538
$ws.Cells.Item(2, 10).Value = 0.1390105253164009
$ws.Cells.Item(2, 11).Value = 0.4069373316475264
$ws.Cells.Item(2, 15).Value = 3.235120941222391

$ws.Cells.Item(3, 2).Value = 0.3354796680974914
$ws.Cells.Item(3, 3).Value = 0.05160864674414256
$ws.Cells.Item(3, 4).Value = 0.6430604655866148
$ws.Cells.Item(3, 5).Value = 0.261602815138879
$ws.Cells.Item(3, 7).Value = 0.7628159505464822
$ws.Cells.Item(3, 8).Value = 0.8608410223957037
$ws.Cells.Item(3, 9).Value = 0.6791021183086272
$ws.Cells.Item(3, 10).Value = 0.1355803821831572
$ws.Cells.Item(3, 11).Value = 0.3615691372934862
$ws.Cells.Item(3, 15).Value = 3.258303670722526

$ws.Cells.Item(4, 2).Value = 0.310730474835907
$ws.Cells.Item(4, 3).Value = 0.04716611989339015
$ws.Cells.Item(4, 4).Value = 0.6368812712590284
$ws.Cells.Item(4, 5).Value = 0.2585573361517319
$ws.Cells.Item(4, 7).Value = 0.7662127791840945
$ws.Cells.Item(4, 8).Value = 0.8651462665618084
$ws.Cells.Item(4, 9).Value = 0.6839504632040629
$ws.Cells.Item(4, 10).Value = 0.1335514430322533
$ws.Cells.Item(4, 11).Value = 0.3336912573301447
$ws.Cells.Item(4, 15).Value = 3.274333397602078

$ws.Cells.Item(5, 2).Value = 0.3006445635429884
$ws.Cells.Item(5, 3).Value = 0.04535299568922824
$ws.Cells.Item(5, 4).Value = 0.6344458353563311
$ws.Cells.Item(5, 5).Value = 0.2573515281882877
$ws.Cells.Item(5, 7).Value = 0.7677197342942534
$ws.Cells.Item(5, 8).Value = 0.866993218043973
$ws.Cells.Item(5, 9).Value = 0.6860174023932828
$ws.Cells.Item(5, 10).Value = 0.1327440543458209
$ws.Cells.Item(5, 11).Value = 0.3223260556349317
$ws.Cells.Item(5, 15).Value = 3.281317013569961

$ws.Cells.Item(6, 2).Value = 0.2989697998208669
$ws.Cells.Item(6, 3).Value = 0.04505176399733557
$ws.Cells.Item(6, 4).Value = 0.63404642872473
$ws.Cells.Item(6, 5).Value = 0.2571534351454687
$ws.Cells.Item(6, 7).Value = 0.7679773722979277
$ws.Cells.Item(6, 8).Value = 0.8673054928777049
$ws.Cells.Item(6, 9).Value = 0.6863661232601999
$ws.Cells.Item(6, 10).Value = 0.1326111613868264
$ws.Cells.Item(6, 11).Value = 0.320438606934772
$ws.Cells.Item(6, 15).Value = 3.282503896389159

$ws.Cells.Item(7, 2).Value = 0.3105944535222704
$ws.Cells.Item(7, 3).Value = 0.04714167850998763
$ws.Cells.Item(7, 4).Value = 0.6368480913195071
$ws.Cells.Item(7, 5).Value = 0.2585409314137124
$ws.Cells.Item(7, 7).Value = 0.7662326057439088
$ws.Cells.Item(7, 8).Value = 0.8651708004853731
$ws.Cells.Item(7, 9).Value = 0.6839779694535437
$ws.Cells.Item(7, 10).Value = 0.1335404756562681
$ws.Cells.Item(7, 11).Value = 0.3335380005069055
$ws.Cells.Item(7, 15).Value = 3.274425753644962

$ws.Cells.Item(8, 2).Value = 0.3618861230415575
$ws.Cells.Item(8, 3).Value = 0.05633954143590358
$ws.Cells.Item(8, 4).Value = 0.6499362681628043
$ws.Cells.Item(8, 5).Value = 0.2649729759241453
$ws.Cells.Item(8, 7).Value = 0.7596109619683133
$ws.Cells.Item(8, 8).Value = 0.8565630776064737
$ws.Cells.Item(8, 9).Value = 0.6742403186285415
$ws.Cells.Item(8, 10).Value = 0.1378117986937681
$ws.Cells.Item(8, 11).Value = 0.3912992709146579
$ws.Cells.Item(8, 15).Value = 3.242741637947944

$ws.Cells.Item(9, 2).Value = 0.4624118025501502
$ws.Cells.Item(9, 3).Value = 0.07428318741399664
$ws.Cells.Item(9, 4).Value = 0.6781998659370458
$ws.Cells.Item(9, 5).Value = 0.2786956786062476
$ws.Cells.Item(9, 7).Value = 0.7505059007826986
$ws.Cells.Item(9, 8).Value = 0.8425996096292323
$ws.Cells.Item(9, 9).Value = 0.6580251738194356
$ws.Cells.Item(9, 10).Value = 0.1468003282414472
$ws.Cells.Item(9, 11).Value = 0.5043722343807815
$ws.Cells.Item(9, 15).Value = 3.194861364887601

$ws.Cells.Item(10, 2).Value = 0.5362088511320167
$ws.Cells.Item(10, 3).Value = 0.08740633325083991
$ws.Cells.Item(10, 4).Value = 0.700545353938594
$ws.Cells.Item(10, 5).Value = 0.2894539238432117
$ws.Cells.Item(10, 7).Value = 0.7461890025880251
$ws.Cells.Item(10, 8).Value = 0.8341160146950699
$ws.Cells.Item(10, 9).Value = 0.6478699235451835
$ws.Cells.Item(10, 10).Value = 0.1537785785617984
$ws.Cells.Item(10, 11).Value = 0.5873014490630908
$ws.Cells.Item(10, 15).Value = 3.168382915353362

$ws.Cells.Item(11, 2).Value = 0.569763506694926
$ws.Cells.Item(11, 3).Value = 0.0933627681798157
$ws.Cells.Item(11, 4).Value = 0.711053492435866
$ws.Cells.Item(11, 5).Value = 0.2944950012080696
$ws.Cells.Item(11, 7).Value = 0.7447419307005703
$ws.Cells.Item(11, 8).Value = 0.8306416187794383
$ws.Cells.Item(11, 9).Value = 0.6436322298214208
$ws.Cells.Item(11, 10).Value = 0.1570347245516928
$ws.Cells.Item(11, 11).Value = 0.6249916142176346
$ws.Cells.Item(11, 15).Value = 3.158228605317902

$ws.Cells.Item(12, 2).Value = 0.5824669272318772
$ws.Cells.Item(12, 3).Value = 0.0956163177064866
$ws.Cells.Item(12, 4).Value = 0.7150818713779472
$ws.Cells.Item(12, 5).Value = 0.2964250504884305
$ws.Cells.Item(12, 7).Value = 0.7442683835499793
$ws.Cells.Item(12, 8).Value = 0.8293812568708177
$ws.Cells.Item(12, 9).Value = 0.6420824984156823
$ws.Cells.Item(12, 10).Value = 0.1582794932083971
$ws.Cells.Item(12, 11).Value = 0.6392583100661113
$ws.Cells.Item(12, 15).Value = 3.154655519383113

$ws.Cells.Item(13, 2).Value = 0.5797311628316208
$ws.Cells.Item(13, 3).Value = 0.09513106711469277
$ws.Cells.Item(13, 4).Value = 0.7142121033348019
$ws.Cells.Item(13, 5).Value = 0.2960084420908444
$ws.Cells.Item(13, 7).Value = 0.7443670582206039
$ws.Cells.Item(13, 8).Value = 0.8296502384961002
$ws.Cells.Item(13, 9).Value = 0.6424138140508298
$ws.Cells.Item(13, 10).Value = 0.1580108882044868
$ws.Cells.Item(13, 11).Value = 0.6361859896439626
$ws.Cells.Item(13, 15).Value = 3.155412940681856

$ws.Cells.Item(14, 2).Value = 0.5708086888566584
$ws.Cells.Item(14, 3).Value = 0.09354821019195469
$ws.Cells.Item(14, 4).Value = 0.7113839249591365
$ws.Cells.Item(14, 5).Value = 0.2946533648806948
$ws.Cells.Item(14, 7).Value = 0.7447014793345375
$ws.Cells.Item(14, 8).Value = 0.8305368195492662
$ws.Cells.Item(14, 9).Value = 0.6435036300195627
$ws.Cells.Item(14, 10).Value = 0.1571368972650333
$ws.Cells.Item(14, 11).Value = 0.6261654634562035
$ws.Cells.Item(14, 15).Value = 3.157929190186508

$ws.Cells.Item(15, 2).Value = 0.5653430045993559
$ws.Cells.Item(15, 3).Value = 0.09257839809225743
$ws.Cells.Item(15, 4).Value = 0.7096579833980741
$ws.Cells.Item(15, 5).Value = 0.293826087673466
$ws.Cells.Item(15, 7).Value = 0.7449160183097092
$ws.Cells.Item(15, 8).Value = 0.8310870792203389
$ws.Cells.Item(15, 9).Value = 0.6441783371168199
$ws.Cells.Item(15, 10).Value = 0.1566030806950351
$ws.Cells.Item(15, 11).Value = 0.6200268304409349
$ws.Cells.Item(15, 15).Value = 3.159505911837527

$ws.Cells.Item(16, 2).Value = 0.5340155938061173
$ws.Cells.Item(16, 3).Value = 0.08701678873549668
$ws.Cells.Item(16, 4).Value = 0.6998655122677064
$ws.Cells.Item(16, 5).Value = 0.2891274333139791
$ws.Cells.Item(16, 7).Value = 0.7462939803667155
$ws.Cells.Item(16, 8).Value = 0.8343508165737177
$ws.Cells.Item(16, 9).Value = 0.6481545585987938
$ws.Cells.Item(16, 10).Value = 0.1535674248348045
$ws.Cells.Item(16, 11).Value = 0.5848375442858469
$ws.Cells.Item(16, 15).Value = 3.169084591874963

$ws.Cells.Item(17, 2).Value = 0.5147926473535733
$ws.Cells.Item(17, 3).Value = 0.08360142857102915
$ws.Cells.Item(17, 4).Value = 0.693945903573649
$ws.Cells.Item(17, 5).Value = 0.2862826009600283
$ws.Cells.Item(17, 7).Value = 0.7472717461647704
$ws.Cells.Item(17, 8).Value = 0.8364515630062925
$ws.Cells.Item(17, 9).Value = 0.6506917208353471
$ws.Cells.Item(17, 10).Value = 0.1517260662463542
$ws.Cells.Item(17, 11).Value = 0.5632406186834658
$ws.Cells.Item(17, 15).Value = 3.175445260046388

$ws.Cells.Item(18, 2).Value = 0.5037346533906941
$ws.Cells.Item(18, 3).Value = 0.08163575333061601
$ws.Cells.Item(18, 4).Value = 0.6905734014506777
$ws.Cells.Item(18, 5).Value = 0.2846601771544854
$ws.Cells.Item(18, 7).Value = 0.7478827585989904
$ws.Cells.Item(18, 8).Value = 0.8376960804887688
$ws.Cells.Item(18, 9).Value = 0.6521869795048438
$ws.Cells.Item(18, 10).Value = 0.1506746566880111
$ws.Cells.Item(18, 11).Value = 0.5508154123265285
$ws.Cells.Item(18, 15).Value = 3.179281712242528

$ws.Cells.Item(19, 2).Value = 0.4999903766273519
$ws.Cells.Item(19, 3).Value = 0.0809699978681806
$ws.Cells.Item(19, 4).Value = 0.6894370825580154
$ws.Cells.Item(19, 5).Value = 0.2841132323397133
$ws.Cells.Item(19, 7).Value = 0.7480979848728282
$ws.Cells.Item(19, 8).Value = 0.8381236748645193
$ws.Cells.Item(19, 9).Value = 0.6526994209398751
$ws.Cells.Item(19, 10).Value = 0.1503199888171594
$ws.Cells.Item(19, 11).Value = 0.5466079199275953
$ws.Cells.Item(19, 15).Value = 3.180611227645215

$ws.Cells.Item(20, 2).Value = 0.5168391193734294
$ws.Cells.Item(20, 3).Value = 0.0839651296534214
$ws.Cells.Item(20, 4).Value = 0.6945727142106932
$ws.Cells.Item(20, 5).Value = 0.2865840053109281
$ws.Cells.Item(20, 7).Value = 0.7471626276576444
$ws.Cells.Item(20, 8).Value = 0.8362241861066622
$ws.Cells.Item(20, 9).Value = 0.6504179144313653
$ws.Cells.Item(20, 10).Value = 0.1519212860896317
$ws.Cells.Item(20, 11).Value = 0.5655399876342528
$ws.Cells.Item(20, 15).Value = 3.17474973579877

$ws.Cells.Item(21, 2).Value = 0.5734295227099153
$ws.Cells.Item(21, 3).Value = 0.09401318936772896
$ws.Cells.Item(21, 4).Value = 0.712213296164748
$ws.Cells.Item(21, 5).Value = 0.2950508115820227
$ws.Cells.Item(21, 7).Value = 0.7446012308053866
$ws.Cells.Item(21, 8).Value = 0.8302749081064889
$ws.Cells.Item(21, 9).Value = 0.6431820317044448
$ws.Cells.Item(21, 10).Value = 0.1573932911540936
$ws.Cells.Item(21, 11).Value = 0.6291088969975931
$ws.Cells.Item(21, 15).Value = 3.157182719359213

$ws.Cells.Item(22, 2).Value = 0.6103968177155821
$ws.Cells.Item(22, 3).Value = 0.1005683373689408
$ws.Cells.Item(22, 4).Value = 0.7240289922134195
$ws.Cells.Item(22, 5).Value = 0.3007073382543766
$ws.Cells.Item(22, 7).Value = 0.7433610885731952
$ws.Cells.Item(22, 8).Value = 0.8267091273886535
$ws.Cells.Item(22, 9).Value = 0.6387735242383599
$ws.Cells.Item(22, 10).Value = 0.161037968238162
$ws.Cells.Item(22, 11).Value = 0.6706209805712433
$ws.Cells.Item(22, 15).Value = 3.147287950124621

$ws.Cells.Item(23, 2).Value = 0.5906685362753308
$ws.Cells.Item(23, 3).Value = 0.09707084890601436
$ws.Cells.Item(23, 4).Value = 0.7176965650274951
$ws.Cells.Item(23, 5).Value = 0.2976771074120421
$ws.Cells.Item(23, 7).Value = 0.7439832345256292
$ws.Cells.Item(23, 8).Value = 0.8285827578016551
$ws.Cells.Item(23, 9).Value = 0.6410970755879255
$ws.Cells.Item(23, 10).Value = 0.1590864804779812
$ws.Cells.Item(23, 11).Value = 0.6484685483081591
$ws.Cells.Item(23, 15).Value = 3.152423754391407

$ws.Cells.Item(24, 2).Value = 0.5159139294622719
$ws.Cells.Item(24, 3).Value = 0.08380070705001685
$ws.Cells.Item(24, 4).Value = 0.6942892372981078
$ws.Cells.Item(24, 5).Value = 0.2864476995704379
$ws.Cells.Item(24, 7).Value = 0.7472118078791397
$ws.Cells.Item(24, 8).Value = 0.836326868662141
$ws.Cells.Item(24, 9).Value = 0.6505415882714445
$ws.Cells.Item(24, 10).Value = 0.1518330047335894
$ws.Cells.Item(24, 11).Value = 0.5645004704302892
$ws.Cells.Item(24, 15).Value = 3.175063622819579

$ws.Cells.Item(25, 2).Value = 0.4352255861202252
$ws.Cells.Item(25, 3).Value = 0.06943927557834684
$ws.Cells.Item(25, 4).Value = 0.6702760621325865
$ws.Cells.Item(25, 5).Value = 0.2748645655381239
$ws.Cells.Item(25, 7).Value = 0.7525529318504312
$ws.Cells.Item(25, 8).Value = 0.8460651488973383
$ws.Cells.Item(25, 9).Value = 0.6621032239570113
$ws.Cells.Item(25, 10).Value = 0.1443030193720745
$ws.Cells.Item(25, 11).Value = 0.4738067166613575
$ws.Cells.Item(25, 15).Value = 3.206287348848264
